# Generate Report for Handoff
#
# The "Ready for handoff" files (rows 4-7 on the per-locale sheets) had
# priority "low"; handing them off bumps their priority to "ht" and
# records the new handoff timestamp in column H.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-26 18:32:37"

    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-26 18:32:41"
}
